$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename sheets (the headline change - "data" preparation pipeline naming)
#    Sheet1 (3) -> percentages
#    Sheet1 (2) -> cleaned
#    Sheet1      -> original
#    ("final" keeps its name)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Sheet1 (3)").Name = "percentages"
$wb.Worksheets.Item("Sheet1 (2)").Name = "cleaned"
$wb.Worksheets.Item("Sheet1").Name = "original"

# ---------------------------------------------------------------------------
# 2) On the "percentages" sheet, the view had scrolled down (topLeftCell
#    A39); reset the scroll position back to the top-left (A1) while leaving
#    the actual selection (I68) untouched.
# ---------------------------------------------------------------------------
$percentages = $wb.Worksheets.Item("percentages")
$percentages.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# ---------------------------------------------------------------------------
# 3) On the "final" sheet, column A carried a redundant/empty direct column
#    format (no actual alignment set). Remove that column-level formatting
#    while preserving every individual cell's existing formatting exactly
#    (header cell A1, header cells B1:D1, and the body cells A2:A194).
#
#    We capture each distinct existing format onto scratch cells first, wipe
#    the column formatting, and then restore the captured formats - this
#    keeps cell styles identical instead of rebuilding them property by
#    property (which would create new, slightly different style entries).
# ---------------------------------------------------------------------------
$final = $wb.Worksheets.Item("final")

$final.Range("A1").Copy() | Out-Null
$final.Range("ZZ1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$final.Range("A2").Copy() | Out-Null
$final.Range("ZZ2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$final.Columns.Item(1).ClearFormats()

$final.Range("ZZ1").Copy() | Out-Null
$final.Range("A1").PasteSpecial(-4122) | Out-Null

$final.Range("ZZ2").Copy() | Out-Null
$final.Range("A2:A194").PasteSpecial(-4122) | Out-Null

$final.Range("ZZ1:ZZ2").Clear()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Restore "final" as the active / selected sheet (it was the active sheet
#    before we activated "percentages" above).
# ---------------------------------------------------------------------------
$final.Activate()
$final.Range("D2").Select()
